# Karnaugh-Map workbook update: add a 4th input variable (D) to the truth
# table on "Plan1" and rebuild the Karnaugh map on "Planilha1" accordingly.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Plan1")
$ws2 = $wb.Worksheets.Item("Planilha1")

# ---------------------------------------------------------------------
# 1) New / relocated header text.
#    Order below matters: it reproduces the exact order in which Excel
#    appends brand-new entries to the shared-string table.
# ---------------------------------------------------------------------

# Move the old "S" header out of D1 into the new E1 before D1 is reused.
$ws1.Range("E1").Value = "S"
$ws1.Range("A1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)

# Planilha1 header / row labels.
$ws2.Range("A1").Value = "V"

$ws2.Range("A2").Formula = "'""-D-A"""
$ws2.Range("A3").Value = '"-DA"'
$ws2.Range("A4").Value = '"DA"'
$ws2.Range("A5").Value = '"D-A"'

# New "D" column header on Plan1 (reuses the D1 cell that used to hold "S").
$ws1.Range("D1").Value = "D"

$ws2.Range("C1").Value = '"-BC"'
$ws2.Range("B1").Value = '"-B-C"'
$ws2.Range("D1").Value = '"BC"'
$ws2.Range("E1").Value = '"B-C"'

# ---------------------------------------------------------------------
# 2) Plan1 truth table: extend to 16 rows (4 inputs: A,B,C,D -> S)
# ---------------------------------------------------------------------

$ws1.Columns.Item(4).ColumnWidth = 8.3

# Column E (S results) for the original 8 rows, matching existing border
# style used by the rest of the table.
$ws1.Range("A2").Copy()
$ws1.Range("E2:E9").PasteSpecial(-4122)

$truth = @(
    @(0,0,0,0,1),
    @(0,0,1,0,1),
    @(0,1,0,0,0),
    @(0,1,1,0,1),
    @(1,0,0,0,1),
    @(1,0,1,0,1),
    @(1,1,0,0,0),
    @(1,1,1,0,1),
    @(0,0,0,1,0),
    @(0,0,1,1,0),
    @(0,1,0,1,0),
    @(0,1,1,1,0),
    @(1,0,0,1,0),
    @(1,0,1,1,0),
    @(1,1,0,1,0),
    @(1,1,1,1,0)
)

# Rows 10-17 are brand new; give them the same formatting as the rest of
# the table before filling in their values.
$ws1.Range("A2:E2").Copy()
$ws1.Range("A10:E17").PasteSpecial(-4122)

for ($i = 0; $i -lt $truth.Length; $i++) {
    $r = $i + 2
    $row = $truth[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------
# 3) Planilha1 Karnaugh map body (4 rows x 4 cols of formulas)
# ---------------------------------------------------------------------

# Rows 4 and 5 are brand new; copy the formatting used by row 3 first.
$ws2.Range("A3:E3").Copy()
$ws2.Range("A4:E5").PasteSpecial(-4122)

$ws2.Range("B2").Formula = "=(Plan1!E2)"
$ws2.Range("C2").Formula = "=(Plan1!E3)"
$ws2.Range("D2").Formula = "=(Plan1!E5)"
$ws2.Range("E2").Formula = "=(Plan1!E4)"

$ws2.Range("B3").Formula = "=(Plan1!E6)"
$ws2.Range("C3").Formula = "=(Plan1!E7)"
$ws2.Range("D3").Formula = "=(Plan1!E9)"
$ws2.Range("E3").Formula = "=(Plan1!E8)"

$ws2.Range("B4").Formula = "=(Plan1!E14)"
$ws2.Range("C4").Formula = "=(Plan1!E15)"
$ws2.Range("D4").Formula = "=(Plan1!E17)"
$ws2.Range("E4").Formula = "=(Plan1!E16)"

$ws2.Range("B5").Formula = "=(Plan1!E10)"
$ws2.Range("C5").Formula = "=(Plan1!E11)"
$ws2.Range("D5").Formula = "=(Plan1!E13)"
$ws2.Range("E5").Formula = "=(Plan1!E12)"

# ---------------------------------------------------------------------
# 4) Selections (Plan1 selection first, Planilha1 last so Planilha1
#    remains the active/visible tab, matching the saved workbook).
# ---------------------------------------------------------------------

$ws1.Range("E10").Select()
$ws2.Range("D3").Select()
